$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The test-case step blocks (Steps / Expected Results) for TC2, TC3 and TC4
# are rotated: TC2 now shows what used to be TC3's content, TC3 now shows
# what used to be TC4's content, and TC4 now shows what used to be TC2's
# content. The "TC2"/"TC3"/"TC4" labels themselves (B14, B21, B28) stay put.

$ws.Range("B18").Value = "Beneficiário Clica em detalhar diária."
$ws.Range("D18").Value = "SYSTEM Apresenta a tela de Detalhar Diárias"

$ws.Range("B25").Value = "Beneficiário Clica em cancelar diária."
$ws.Range("D25").Value = "SYSTEM Apresenta a tela de Cancelar Solicitação de Diária"

$ws.Range("B32").Value = "Beneficiário Clica em analisar prestação de contas."
$ws.Range("D32").Value = "SYSTEM Apresenta a tela de Analisar Prestação de Contas"
